# "before revise topopoints to incorporate data controls and properties from controlpts"
#
# The ΔXY "Slope" (J/L) column is moved two columns to the right (M -> O) and a
# new ΔXYZ column (3-D distance = SQRT(ΔXY^2 + ΔZ^2), i.e. SQRT(L*L+J*J)) is
# inserted in its place, with a narrow spacer column left between the two data
# columns - mirroring the existing J/K/L spacer pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new, blank columns at M:N. This pushes the existing "Slope"
# column (M) and the blank style-only columns after it (N,O,P) two places to
# the right, becoming O,P,Q,R - and keeps every formula reference (J/L, etc.)
# correctly repointed by Excel automatically.
$ws.Columns("M:N").Insert()

# Header for the new distance column.
$ws.Range("M4").Value = "ΔXYZ"

# Row 6 is a plain (non-shared) formula in the source file.
$ws.Range("M6").Formula = "=SQRT(L6*L6+J6*J6)"

# Rows 7-15 form the shared-formula block (matches the existing L7:L19 /
# J-column pattern); rows 16-19 are left blank, same as the source (those
# rows only have the moved "Slope" formula, now living in column O).
$ws.Range("M7:M15").Formula = "=SQRT(L7*L7+J7*J7)"

# Narrow spacer column N, matching the existing spacer-column styling (e.g.
# column K before it).
$ws.Columns("N:N").ColumnWidth = 3.77734375

# Restore/relocate the active selection to M12.
$ws.Range("M12").Select()
